$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 64: Forged from the Void / Void Glue
$ws.Range("H64").Value = 4500
$ws.Range("I64").Value = 4000
$ws.Range("J64").Value = 5500
$ws.Range("K64").Value = 4000
$ws.Range("L64").Value = 5500
$ws.Range("M64").Value = -3752
$ws.Range("N64").Value = -5996

# ALC row 67: Dodging the Draft (L) / Void Glue
$ws.Range("H67").Value = 4500
$ws.Range("I67").Value = 4000
$ws.Range("J67").Value = 5500
$ws.Range("K67").Value = 4000
$ws.Range("L67").Value = 5500
$ws.Range("M67").Value = -3142
$ws.Range("N67").Value = -7216

# ALC row 76: Warding Off Temptation / Enchanted Hardsilver Ink
$ws.Range("H76").Value = 5492.8125
$ws.Range("I76").Value = 5493.1333
$ws.Range("J76").Value = 5488
$ws.Range("K76").Value = 5493.1333
$ws.Range("L76").Value = 5488
$ws.Range("M76").Value = -5178.1333
$ws.Range("N76").Value = -6118

# ALC row 79: The Garden of Arcane Delights (L) / Enchanted Hardsilver Ink
$ws.Range("H79").Value = 5492.8125
$ws.Range("I79").Value = 5493.1333
$ws.Range("J79").Value = 5488
$ws.Range("K79").Value = 5493.1333
$ws.Range("L79").Value = 5488
$ws.Range("M79").Value = -4401.1333
$ws.Range("N79").Value = -7672

# ALC row 116: Growing Up / Growth Formula Kappa
$ws.Range("H116").Value = 3448.2222
$ws.Range("I116").Value = 2838
$ws.Range("J116").Value = 4668.6665
$ws.Range("K116").Value = 2838
$ws.Range("L116").Value = 4668.6665
$ws.Range("M116").Value = 604
$ws.Range("N116").Value = -11552.6665

# ALC row 132: Fast-forwarding Flora / Growth Formula Lambda
$ws.Range("H132").Value = 971.7143
$ws.Range("I132").Value = 954.1539
$ws.Range("J132").Value = 1200
$ws.Range("K132").Value = 2862.4617
$ws.Range("L132").Value = 3600
$ws.Range("M132").Value = -332.4616999999998
$ws.Range("N132").Value = -8660

$ws = $wb.Worksheets.Item("ARM")
# ARM row 18: Still the Best / Brass Alembic
$ws.Range("H18").Value = 200
$ws.Range("I18").Value = 200
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 200
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = 122

# ARM row 102: Smells of Rich Tama-hagane / Tama-hagane Ingot
$ws.Range("H102").Value = 1622.4286
$ws.Range("I102").Value = 1728
$ws.Range("J102").Value = 250
$ws.Range("K102").Value = 1728
$ws.Range("L102").Value = 250
$ws.Range("M102").Value = -106
$ws.Range("N102").Value = -3494

# ARM row 110: Scheduled Maintenance / Deepgold Ingot
$ws.Range("H110").Value = 1105.1111
$ws.Range("I110").Value = 1105.1111
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 1105.1111
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 939.8888999999999

# ARM row 128: Heading toward Bankruptcy / Manganese Helm of the Falling Dragon
$ws.Range("H128").Value = 49999.5
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 49999.5
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 49999.5
$ws.Range("N128").Value = -59959.5

$ws = $wb.Worksheets.Item("BSM")
# BSM row 86: Through Thick and Thin / Adamantite Nugget
$ws.Range("H86").Value = 2284.611
$ws.Range("I86").Value = 839.7143
$ws.Range("J86").Value = 3204.0908
$ws.Range("K86").Value = 839.7143
$ws.Range("L86").Value = 3204.0908
$ws.Range("M86").Value = 283.2857
$ws.Range("N86").Value = -5450.0908

# BSM row 89: Piercing Eyes Deserve Piercing Shafts (L) / Adamantite Nugget
$ws.Range("H89").Value = 2284.611
$ws.Range("I89").Value = 839.7143
$ws.Range("J89").Value = 3204.0908
$ws.Range("K89").Value = 4198.5715
$ws.Range("L89").Value = 16020.454
$ws.Range("M89").Value = 1417.4285
$ws.Range("N89").Value = -27252.454

# BSM row 99: Meddle in Metal / Oroshigane Ingot
$ws.Range("H99").Value = 3315.7693
$ws.Range("I99").Value = 2065.8333
$ws.Range("J99").Value = 4387.143
$ws.Range("K99").Value = 2065.8333
$ws.Range("L99").Value = 4387.143
$ws.Range("M99").Value = -567.8332999999998
$ws.Range("N99").Value = -7383.143

# BSM row 105: Ingot to Wing It / Molybdenum Ingot
$ws.Range("H105").Value = 3257.7693
$ws.Range("I105").Value = 2737.1
$ws.Range("J105").Value = 4993.3335
$ws.Range("K105").Value = 2737.1
$ws.Range("L105").Value = 4993.3335
$ws.Range("M105").Value = -990.0999999999999
$ws.Range("N105").Value = -8487.333500000001

# BSM row 107: The Gold Experience / Deepgold Nugget
$ws.Range("H107").Value = 4839
$ws.Range("I107").Value = 4521.6
$ws.Range("J107").Value = 8013
$ws.Range("K107").Value = 4521.6
$ws.Range("L107").Value = 8013
$ws.Range("M107").Value = -2601.6
$ws.Range("N107").Value = -11853

$ws = $wb.Worksheets.Item("CRP")
# CRP row 58: You Do the Heavy Lifting / Mahogany Lumber
$ws.Range("H58").Value = 5000
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 5000
$ws.Range("K58").Value = 0
$ws.Range("N58").Value = -5406
$ws.Range("L58").ClearContents()
$ws.Range("M58").ClearContents()

# CRP row 105: Zelkova, My Love / Zelkova Lumber
$ws.Range("H105").Value = 2862.5862
$ws.Range("I105").Value = 1986.1428
$ws.Range("J105").Value = 3680.6
$ws.Range("K105").Value = 1986.1428
$ws.Range("L105").Value = 3680.6
$ws.Range("M105").Value = -239.1428000000001
$ws.Range("N105").Value = -7174.6

# CRP row 107: Built to Last / White Oak Lumber
$ws.Range("H107").Value = 698.6875
$ws.Range("I107").Value = 574.5
$ws.Range("J107").Value = 822.875
$ws.Range("K107").Value = 574.5
$ws.Range("L107").Value = 822.875
$ws.Range("M107").Value = 1345.5
$ws.Range("N107").Value = -4662.875

# CRP row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Range("H132").Value = 1819.1765
$ws.Range("I132").Value = 1433.5555
$ws.Range("J132").Value = 2253
$ws.Range("K132").Value = 4300.666499999999
$ws.Range("L132").Value = 6759
$ws.Range("M132").Value = -1770.666499999999
$ws.Range("N132").Value = -11819

# CRP row 136: Turali Quality / Dark Mahogany Lumber
$ws.Range("H136").Value = 5000
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 0
$ws.Range("N136").Value = -20100
$ws.Range("L136").ClearContents()
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# CUL row 5: What a Sap / Maple Syrup
$ws.Range("H5").Value = 383.75
$ws.Range("I5").Value = 416.66666
$ws.Range("J5").Value = 285
$ws.Range("K5").Value = 1249.99998
$ws.Range("L5").Value = 855
$ws.Range("M5").Value = -1137.99998
$ws.Range("N5").Value = -1079

# CUL row 14: Keep Your Powder Dry / Kukuru Powder
$ws.Range("H14").Value = 1948.25
$ws.Range("I14").Value = 1948.25
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 5844.75
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -5671.75

# CUL row 92: Oh No Udon / Gyr Abanian Flour
$ws.Range("H92").Value = 879.2
$ws.Range("I92").Value = 849.25
$ws.Range("J92").Value = 999
$ws.Range("K92").Value = 2547.75
$ws.Range("L92").Value = 2997
$ws.Range("M92").Value = -1299.75
$ws.Range("N92").Value = -5493

# CUL row 94: All You Can Stomach / Baklava
$ws.Range("H94").Value = 15641.429
$ws.Range("I94").Value = 4995
$ws.Range("J94").Value = 19900
$ws.Range("K94").Value = 14985
$ws.Range("L94").Value = 59700
$ws.Range("M94").Value = -14309
$ws.Range("N94").Value = -61052

# CUL row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Range("H131").Value = 2630.875
$ws.Range("I131").Value = 1978.8
$ws.Range("J131").Value = 2927.2727
$ws.Range("K131").Value = 5936.4
$ws.Range("L131").Value = 8781.8181
$ws.Range("M131").Value = -896.3999999999996
$ws.Range("N131").Value = -18861.8181

# CUL row 135: Not-so-secret Ingredient / Royal Maple Syrup
$ws.Range("H135").Value = 383.75
$ws.Range("I135").Value = 416.66666
$ws.Range("J135").Value = 285
$ws.Range("K135").Value = 3749.99994
$ws.Range("L135").Value = 2565
$ws.Range("M135").Value = -1214.99994
$ws.Range("N135").Value = -7635

$ws = $wb.Worksheets.Item("GSM")
# GSM row 34: All Booked Up / Silver Magnifiers
$ws.Range("H34").Value = 36086.5
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 36086.5
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 36086.5
$ws.Range("N34").Value = -36622.5

# GSM row 76: The Monuments Mages / Hardsilver Magnifiers of Casting
$ws.Range("H76").Value = 36086.5
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 36086.5
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 36086.5
$ws.Range("N76").Value = -36716.5

# GSM row 79: Deal with It (L) / Hardsilver Magnifiers of Casting
$ws.Range("H79").Value = 36086.5
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 36086.5
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 36086.5
$ws.Range("N79").Value = -38270.5

# GSM row 126: Gold Rush Order / Phrygian Gold Ingot
$ws.Range("H126").Value = 4703
$ws.Range("I126").Value = 4484.2
$ws.Range("J126").Value = 5250
$ws.Range("K126").Value = 13452.6
$ws.Range("L126").Value = 15750
$ws.Range("M126").Value = -10982.6
$ws.Range("N126").Value = -20690

$ws = $wb.Worksheets.Item("LTW")
# LTW row 5: These Boots Are Made for Wailing / Leather Duckbills of Gathering
$ws.Range("H5").Value = 39999
$ws.Range("I5").Value = 39999
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 39999
$ws.Range("L5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("N5").ClearContents()

# LTW row 40: Best Served Toad / Toad Leather
$ws.Range("H40").Value = 3255.2222
$ws.Range("I40").Value = 2882.1875
$ws.Range("J40").Value = 3797.818
$ws.Range("K40").Value = 2882.1875
$ws.Range("L40").Value = 3797.818
$ws.Range("M40").Value = -2746.1875
$ws.Range("N40").Value = -4069.818

# LTW row 132: Tenets of Tanning / Silver Lobo Leather
$ws.Range("H132").Value = 3581.3333
$ws.Range("I132").Value = 2541
$ws.Range("J132").Value = 4621.6665
$ws.Range("K132").Value = 7623
$ws.Range("L132").Value = 13864.9995
$ws.Range("M132").Value = -5093
$ws.Range("N132").Value = -18924.9995

$ws = $wb.Worksheets.Item("WVR")
# WVR row 68: What Not to Wear / Holy Rainbow Shirt of Striking
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").ClearContents()
$ws.Range("N68").ClearContents()

# WVR row 71: Appeal of Foreign Apparel (L) / Holy Rainbow Shirt of Striking
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").ClearContents()
$ws.Range("N71").ClearContents()

# WVR row 81: Where the Dragonflies, the Net Catches / Crawler Silk
$ws.Range("H81").Value = 717032.6
$ws.Range("I81").Value = 1945.6
$ws.Range("J81").Value = 2504750.2
$ws.Range("K81").Value = 3891.2
$ws.Range("L81").Value = 5009500.4
$ws.Range("M81").Value = -2830.2
$ws.Range("N81").Value = -5011622.4

# WVR row 84: To Kill a Dragon on Nameday (L) / Crawler Silk
$ws.Range("H84").Value = 717032.6
$ws.Range("I84").Value = 1945.6
$ws.Range("J84").Value = 2504750.2
$ws.Range("K84").Value = 19456
$ws.Range("L84").Value = 25047502
$ws.Range("M84").Value = -14152
$ws.Range("N84").Value = -25058110

# WVR row 100: Of Great Import / Kudzu Thread
$ws.Range("H100").Value = 5264035.5
$ws.Range("I100").Value = 7693052.5
$ws.Range("J100").Value = 1165.6666
$ws.Range("K100").Value = 15386105
$ws.Range("L100").Value = 2331.3332
$ws.Range("M100").Value = -15385564
$ws.Range("N100").Value = -3413.3332

# WVR row 107: Flax Wax / Bright Linen Yarn
$ws.Range("H107").Value = 264.14285
$ws.Range("I107").Value = 208.16667
$ws.Range("J107").Value = 600
$ws.Range("K107").Value = 624.50001
$ws.Range("L107").Value = 1800
$ws.Range("M107").Value = 1295.49999
$ws.Range("N107").Value = -5640

# WVR row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 3342.5715
$ws.Range("I132").Value = 3233
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 9699
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -7169
$ws.Range("N132").Value = -17060
